# Insert a new "Test paragraph" clause immediately after "Section 1",
# at the same outline level (numId=10, ilvl=0) and style ("Normal") as
# the clause that already follows it, pushing the existing clauses
# down by one.
#
# Strategy: locate the first existing paragraph whose text starts with
# "Test paragraph" (our insertion anchor) and call
# Range.InsertParagraphBefore() on it. Word clones the anchor
# paragraph's own pPr (pStyle + numPr) onto the freshly-created empty
# paragraph, so the new clause automatically gets the correct numbering
# (numId=10 / ilvl=0) without spawning a brand-new list instance, the
# way re-applying/ApplyListTemplateWithLevel would.
#
# Note: once InsertParagraphBefore() runs, the COM paragraph object we
# called it on reseats to the newly-created (blank) paragraph rather
# than the original text -- so we resolve the anchor by its numeric
# Paragraphs() index both before and after the insert, instead of
# holding onto the paragraph reference across the call.

$d = $word.ActiveDocument

$anchorIndex = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Test paragraph*") {
        $anchorIndex = $p.Index
        break
    }
}

$anchorParagraph = $d.Paragraphs($anchorIndex)
$anchorParagraph.Range.InsertParagraphBefore()

# The blank paragraph InsertParagraphBefore() created now occupies
# $anchorIndex (it inherited the anchor's pPr/numPr/pStyle); the
# original "Test paragraph" text got pushed down to $anchorIndex + 1.
$newPara = $d.Paragraphs($anchorIndex)
$newPara.Range.Text = "Test paragraph"
